$d = $word.ActiveDocument

# Locate the paragraph "This picture shows the code I alternated." and
# collapse the found range to its end (i.e. right after that sentence).
$rng = $d.Content
$found = $rng.Find.Execute("This picture shows the code I alternated.", `
                            $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the anchor paragraph 'This picture shows the code I alternated.'"
}
$rng.Collapse(0)

# Insert a new blank paragraph right after "This picture shows the code I
# alternated." (this becomes the new paragraph 28, pushing the document's
# existing blank paragraph further down).
$rng.InsertParagraphAfter()
$rng = $d.Range($rng.End + 1, $rng.End + 1)

# Insert a paragraph with "This is the link to the code."
$rng.InsertParagraphAfter()
$rng = $d.Range($rng.End + 1, $rng.End + 1)
$rng.InsertBefore("This is the link to the code.")
$rng = $d.Range($rng.End, $rng.End)

# Insert a paragraph with the link to the Arduino code on GitHub.
$rng.InsertParagraphAfter()
$rng = $d.Range($rng.End + 1, $rng.End + 1)
$rng.InsertBefore("https://github.com/brigham4210/Fan_Arduino/blob/main/FAN/FAN.ino")

Write-Output "Inserted link paragraphs after the code-picture caption."
